$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume-change (E) columns to refreshed values.
# Numeric-looking price strings (e.g. "586.51") must stay plain text cells
# (as in the source data, which uses "." as a thousands separator, e.g.
# "63.249.96"), so we briefly force a Text number format before assigning
# the value and then restore the default "Normal" style.

$ws.Range("D2").Value = "63.249.96"
$ws.Range("E2").Value = "  +5.00%  "

$ws.Range("D3").Value = "2.708.43"
$ws.Range("E3").Value = "  +4.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.60%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.32%  "

$ws.Range("E8").Value = "  +1.32%  "

$ws.Range("D9").Value = "2.738.87"
$ws.Range("E9").Value = "  +5.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.86%  "

$ws.Range("E11").Value = "  +7.36%  "

$ws.Range("E12").Value = "  +4.44%  "

$ws.Range("E13").Value = "  +1.61%  "

$ws.Range("D14").Value = "3.191.87"
$ws.Range("E14").Value = "  +4.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.41%  "

$ws.Range("D16").Value = "63.109.25"
$ws.Range("E16").Value = "  +4.78%  "

$ws.Range("E17").Value = "  +7.24%  "

$ws.Range("D18").Value = "2.727.92"
$ws.Range("E18").Value = "  +4.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.07%  "

$ws.Range("E22").Value = "  +1.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("E24").Value = "  +0.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.79%  "

$ws.Range("E26").Value = "  +3.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.18%  "

$ws.Range("E29").Value = "  +7.56%  "

$ws.Range("E30").Value = "  +6.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.86%  "

$ws.Range("E33").Value = "  +21.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.36%  "

$ws.Range("E38").Value = "  +10.76%  "

$ws.Range("E39").Value = "  +18.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "351.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.22%  "

$ws.Range("E46").Value = "  +7.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "138.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.31%  "

$ws.Range("E48").Value = "  +7.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.640"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.75%  "

$ws.Range("E50").Value = "  +1.69%  "

$ws.Range("E51").Value = "  -0.38%  "

